$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = '*maa://24633 (54.97), *maa://30515 (69.7), ***maa://20792 (11.93), *maa://34787 (72.13), maa://39402 (85.19), ***maa://29083 (29.41)'
$ws.Range("AE2").Value = '***maa://21730 (17.19), maa://25251 (92.11), ***maa://39501 (20.0), *maa://36675 (60.0)'
$ws.Range("AE3").Value = '*maa://21289 (68.42)'
$ws.Range("AE4").Value = '*maa://30062 (63.41), ***maa://26209 (13.04), *maa://39394 (75.0)'
$ws.Range("AA9").Value = 'maa://28711 (87.8), ***maa://22740 (5.88), **maa://27377 (46.15), ***maa://25174 (20.0), *maa://39938 (53.85), maa://40166 (100.0)'
$ws.Range("W11").Value = 'maa://36713 (97.79)'
$ws.Range("K14").Value = 'maa://26245 (96.12), maa://21288 (96.21), maa://36682 (100.0), maa://39841 (92.5)'
$ws.Range("O14").Value = 'maa://23250 (98.46), maa://20107 (87.1), maa://22772 (100.0), **maa://22745 (50.0)'
$ws.Range("AE15").Value = 'maa://21364 (80.61), *maa://22766 (73.0), *maa://36666 (77.05)'
$ws.Range("C18").Value = 'maa://24570 (96.45)'
$ws.Range("AA19").Value = '*maa://30709 (60.27), *maa://36668 (52.17)'
$ws.Range("K20").Value = 'maa://41331 (87.5)'
$ws.Range("AA21").Value = '*maa://21443 (78.77), ***maa://23820 (29.63)'
$ws.Range("K23").Value = 'maa://39756 (92.31), maa://39875 (95.24)'
$ws.Range("O23").Value = 'maa://30587 (91.62), *maa://29748 (75.2), ***maa://29785 (15.15), *maa://37566 (76.47)'
$ws.Range("W24").Value = 'maa://23504 (93.04), maa://29988 (85.86), **maa://22892 (40.43), *maa://25141 (76.86), *maa://36663 (79.63), ***maa://22815 (23.08)'
$ws.Range("C25").Value = 'maa://29753 (95.11)'
$ws.Range("G25").Value = '*maa://29063 (76.74), *maa://25311 (74.73), ***maa://22725 (4.84)'
$ws.Range("AE26").Value = 'maa://30511 (83.87), *maa://29760 (61.54)'
$ws.Range("G27").Value = '**maa://21283 (48.65), maa://34494 (100.0), **maa://36665 (44.44), maa://39601 (85.71)'
$ws.Range("C28").Value = 'maa://24465 (90.32), maa://25725 (81.82)'
$ws.Range("W28").Value = 'maa://39929 (85.8), ***maa://39723 (15.15), *maa://41749 (80.0)'
$ws.Range("AE28").Value = 'maa://36660 (94.05), *maa://36701 (64.0)'
$ws.Range("K29").Value = 'maa://28432 (93.38), *maa://28440 (72.5), maa://31400 (100.0), *maa://28650 (66.67)'
$ws.Range("K31").Value = 'maa://35926 (93.45), maa://36258 (80.26)'
$ws.Range("S32").Value = 'maa://41108 (88.24), maa://41238 (94.12)'
$ws.Range("G46").Value = 'maa://35931 (92.38)'
$ws.Range("G47").Value = 'maa://27410 (95.79), maa://29661 (97.64), maa://28038 (84.62)'
$ws.Range("G55").Value = 'maa://32532 (92.17)'
